$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.256.02"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.830.02"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").Value = "'235.82"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").Value = "'0.6021"
$ws.Range("E6").Value = "  -3.89%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").Value = "'0.07040"
$ws.Range("E8").Value = "  -4.97%  "
$ws.Range("D9").Value = "'0.2790"
$ws.Range("E9").Value = "  -3.44%  "
$ws.Range("D10").Value = "'23.54"
$ws.Range("E10").Value = "  -5.57%  "
$ws.Range("D11").Value = "'0.07647"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "1.831.52"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "'4.787"
$ws.Range("E13").Value = "  -3.63%  "
$ws.Range("D14").Value = "'0.6276"
$ws.Range("E14").Value = "  -6.81%  "
$ws.Range("D15").Value = "'0.000009741"
$ws.Range("E15").Value = "  -4.87%  "
$ws.Range("D16").Value = "'78.99"
$ws.Range("E16").Value = "  -3.33%  "
$ws.Range("D17").Value = "29.227.71"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "'5.824"
$ws.Range("E18").Value = "  -6.09%  "
$ws.Range("D19").Value = "'223.96"
$ws.Range("E19").Value = "  -3.85%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "'11.69"
$ws.Range("E21").Value = "  -5.06%  "
$ws.Range("D22").Value = "'6.999"
$ws.Range("E22").Value = "  -3.98%  "
$ws.Range("D23").Value = "'1.003"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "'156.47"
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("D25").Value = "'7.982"
$ws.Range("E25").Value = "  -5.90%  "
$ws.Range("D26").Value = "'0.1299"
$ws.Range("E26").Value = "  -3.34%  "
$ws.Range("D27").Value = "'16.62"
$ws.Range("E27").Value = "  -3.87%  "
$ws.Range("D28").Value = "'0.06632"
$ws.Range("E28").Value = "  -8.59%  "
$ws.Range("D29").Value = "'1.470"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'1.445"
$ws.Range("E30").Value = "  -1.98%  "
$ws.Range("D31").Value = "'3.843"
$ws.Range("E31").Value = "  -4.51%  "
$ws.Range("D32").Value = "'3.789"
$ws.Range("E32").Value = "  -6.46%  "
$ws.Range("D33").Value = "'1.105"
$ws.Range("E33").Value = "  -2.93%  "
$ws.Range("D34").Value = "'1.722"
$ws.Range("E34").Value = "  -4.89%  "
$ws.Range("D35").Value = "'0.6458"
$ws.Range("E35").Value = "  -7.17%  "
$ws.Range("D36").Value = "'2.550"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").Value = "'2.737"
$ws.Range("E37").Value = "  -2.72%  "
$ws.Range("D38").Value = "1.213.33"
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("D39").Value = "'0.01760"
$ws.Range("E39").Value = "  -4.30%  "
$ws.Range("D40").Value = "'6.532"
$ws.Range("E40").Value = "  -5.33%  "
$ws.Range("D41").Value = "'0.8993"
$ws.Range("E41").Value = "  -5.85%  "
$ws.Range("D42").Value = "'1.005"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").Value = "1.989.50"
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("D44").Value = "'100.38"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").Value = "'62.57"
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("D46").Value = "'0.00000000115"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").Value = "'8.510"
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("D48").Value = "'1.579"
$ws.Range("E48").Value = "  -7.56%  "
$ws.Range("D49").Value = "'0.4558"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").Value = "'0.05504"
$ws.Range("E50").Value = "  -2.67%  "
$ws.Range("D51").Value = "'6.393"
$ws.Range("E51").Value = "  -7.81%  "
